# Fruta / hortaliza, semanal
# Insert two new weekly price records for "Betarraga" (Macroferia Regional de
# Talca) at rows 430-431, pushing the existing data down by two rows
# (old row 430 -> 432, ... old row 508 -> 510), matching the new
# dimension A1:R510.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 430:431 - Excel shifts rows 430..508 down to 432..510
$ws.Range("A430:A431").EntireRow.Insert()

# New row 430
$ws.Range("A430").Value = 5
$ws.Range("B430").Value = 'Macroferia Regional de Talca'
$ws.Range("C430").Value = 'Maule'
$ws.Range("D430").Value = 45015
$ws.Range("E430").Value = 7
$ws.Range("F430").Value = 100114014
$ws.Range("G430").Value = 'Betarraga'
$ws.Range("H430").Value = 'Sin especificar'
$ws.Range("I430").Value = 'Primera'
$ws.Range("J430").Value = 3000
$ws.Range("K430").Value = 600
$ws.Range("L430").Value = 600
$ws.Range("M430").Value = 600
$ws.Range("N430").Value = '$/paquete 5 unidades'
$ws.Range("O430").Value = 'Región del Maule'
$ws.Range("P430").Value = 120
$ws.Range("Q430").Value = 5
$ws.Range("R430").Value = 'Hortaliza'

# New row 431
$ws.Range("A431").Value = 5
$ws.Range("B431").Value = 'Macroferia Regional de Talca'
$ws.Range("C431").Value = 'Maule'
$ws.Range("D431").Value = 45015
$ws.Range("E431").Value = 7
$ws.Range("F431").Value = 100114014
$ws.Range("G431").Value = 'Betarraga'
$ws.Range("H431").Value = 'Sin especificar'
$ws.Range("I431").Value = 'Segunda'
$ws.Range("J431").Value = 2000
$ws.Range("K431").Value = 500
$ws.Range("L431").Value = 500
$ws.Range("M431").Value = 500
$ws.Range("N431").Value = '$/paquete 5 unidades'
$ws.Range("O431").Value = 'Región del Maule'
$ws.Range("P431").Value = 100
$ws.Range("Q431").Value = 5
$ws.Range("R431").Value = 'Hortaliza'
